# Add a new column R (2021 data) to the "3.6.1 Death rate due to road
# traffic injuries" worksheet, mirroring the structure/format of the
# existing column Q (2020 data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values for the new column R, keyed by row number (row 3 is the header
# year row; rows 4-33 are the data rows).
$values = @{
    3  = 2021
    4  = 12.6
    5  = 17.9
    6  = 7.3
    7  = 12.6
    8  = 19.5
    9  = 5.5
    10 = 10.3
    11 = 12.3
    12 = 8.2
    13 = 24.8
    14 = 33.1
    15 = 16.6
    16 = 23.9
    17 = 29.5
    18 = 18.1
    19 = 9.6
    20 = 14.8
    21 = 4.3
    22 = 12.1
    23 = 18.2
    24 = 5.9
    25 = 17.3
    26 = 27.6
    27 = 7.4
    28 = 7.8
    29 = 10.4
    30 = 5.6
    31 = 6.7
    32 = 10.7
    33 = 3
}

# Column Q is column 17, column R is column 18. Copy the formatting of
# each Q-column cell into the corresponding new R-column cell (matching
# number format / font / borders / alignment) before writing the value.
for ($r = 3; $r -le 33; $r++) {
    $src = $ws.Cells.Item($r, 17)
    $dst = $ws.Cells.Item($r, 18)
    $src.Copy($dst)
    $dst.Value = $values[$r]
}

# Match the author's recorded selection after the edit.
$ws.Range("S4").Select()
